# "new coin and fish coords"
# Update several (x,y) coordinate rows, remove rows that no longer have a
# matching record, drop the now-unused trailing index-only rows, and turn
# the data range into a proper AutoFilter table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual coordinate values (B = x, C = y) ---------------
$ws.Range("B3").Value = 106
$ws.Range("C3").Value = 1526

$ws.Range("B5").Value = 828

$ws.Range("B12").Value = 1332

$ws.Range("B21").Value = 1112
$ws.Range("C21").Value = 6664

$ws.Range("B23").Value = 10

$ws.Range("B24").Value = 1141
$ws.Range("C24").Value = 7428

$ws.Range("B26").Value = 1031
$ws.Range("C26").Value = 7946

$ws.Range("B28").Value = 128

$ws.Range("B29").Value = 10
$ws.Range("C29").Value = 8754

# --- Remove rows that correspond to records no longer present ---------
# (cleared, not shifted, so remaining row numbers keep their identity)
$ws.Range("A4:E4").ClearContents()
$ws.Range("A7:E7").ClearContents()
$ws.Range("A10:E10").ClearContents()
$ws.Range("A11:E11").ClearContents()
$ws.Range("A22:E22").ClearContents()

# --- Drop the leftover trailing index-only rows ------------------------
$ws.Range("A33:A36").ClearContents()

# --- Turn the table into a filterable range ----------------------------
$ws.Range("A1:E36").AutoFilter() | Out-Null
$n = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$E`$36")
$n.Visible = $false

# --- Restore the active selection --------------------------------------
$ws.Range("B29").Select() | Out-Null
